$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Formula = "'316.29"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.Formula = "'1.87%"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Formula = "'48.77"
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.Formula = "'10.32%"
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Formula = "'5.299"
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Formula = "'4.17%"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Formula = "'0.07915"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.Formula = "'-0.97%"
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Formula = "'2.65%"
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Formula = "'1.328"
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Formula = "'24.15%"
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Formula = "'1.631"
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.Formula = "'0.16%"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Formula = "'0.1241"
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.Formula = "'-3.61%"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Formula = "'0.1967"
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Formula = "'3.94%"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Formula = "'0.09520"
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Formula = "'3.65%"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Formula = "'0.04559"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Formula = "'8.83%"
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Formula = "'0.1050"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.Formula = "'1.78%"
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Formula = "'0.001312"
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Formula = "'0.54%"
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Formula = "'0.04211"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.Formula = "'0.94%"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Formula = "'0.005907"
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.Formula = "'3.65%"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Formula = "'-0.81%"
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Formula = "'2.471"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Formula = "'2.88%"
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.Formula = "'3.24%"
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Formula = "'8.040"
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.Formula = "'0.78%"
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Formula = "'0.1406"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Formula = "'2.56%"
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Formula = "'0.3073"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.Formula = "'-1.37%"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Formula = "'2.52%"
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Formula = "'0.004186"
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.Formula = "'-3.39%"
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Formula = "'0.0001355"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Formula = "'1.73%"
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Formula = "'0.0003560"
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Formula = "'0.02644"
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.Formula = "'0.20%"
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Formula = "'0.05881"
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.Formula = "'9.25%"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Formula = "'0.01036"
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.Formula = "'85.01%"
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Formula = "'0.008042"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Formula = "'4.26%"
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.Formula = "'3.18%"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Formula = "'0.007559"
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.Formula = "'4.12%"
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Formula = "'0.007938"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.Formula = "'-5.18%"
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Formula = "'0.3164"
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Formula = "'2.32%"
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Formula = "'0.00007041"
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Formula = "'5.16%"
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Formula = "'0.00000000753"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Formula = "'1.75%"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Formula = "'0.05598"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Formula = "'10.25%"
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Formula = "'0.004023"
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Formula = "'1.92%"
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Formula = "'0.00002109"
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Formula = "'1.75%"
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Formula = "'0.0002008"
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.Formula = "'1.75%"
$c.Style = "Normal"
